$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 685.28
$ws.Range("I41").Value = 881.7273
$ws.Range("K41").Value = 881.7273
$ws.Range("M41").Value = -441.7273

# Row 51
$ws.Range("H51").Value = 11726.091
$ws.Range("I51").Value = 3000.5
$ws.Range("J51").Value = 13665.111
$ws.Range("K51").Value = 3000.5
$ws.Range("L51").Value = 13665.111
$ws.Range("M51").Value = -2516.5
$ws.Range("N51").Value = -14633.111

# Row 76
$ws.Range("H76").Value = 3228.6487
$ws.Range("I76").Value = 3040.4443
$ws.Range("K76").Value = 3040.4443
$ws.Range("M76").Value = -2725.4443

# Row 79
$ws.Range("H79").Value = 3228.6487
$ws.Range("I79").Value = 3040.4443
$ws.Range("K79").Value = 3040.4443
$ws.Range("M79").Value = -1948.4443

# Row 138
$ws.Range("H138").Value = 1829.9125
$ws.Range("I138").Value = 1041.7142
$ws.Range("J138").Value = 1997.1061
$ws.Range("K138").Value = 3125.1426
$ws.Range("L138").Value = 5991.3183
$ws.Range("M138").Value = 2014.8574
$ws.Range("N138").Value = -16271.3183


# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4565.4546
$ws.Range("I61").Value = 2357.7778
$ws.Range("J61").Value = 14500
$ws.Range("K61").Value = 2357.7778
$ws.Range("L61").Value = 14500
$ws.Range("M61").Value = -2145.7778
$ws.Range("N61").Value = -14924

# Row 122
$ws.Range("H122").Value = 1073.9166
$ws.Range("I122").Value = 909.6667
$ws.Range("K122").Value = 2729.0001
$ws.Range("M122").Value = -279.0001000000002

# Row 136
$ws.Range("H136").Value = 4565.4546
$ws.Range("I136").Value = 2357.7778
$ws.Range("J136").Value = 14500
$ws.Range("K136").Value = 7073.3334
$ws.Range("L136").Value = 43500
$ws.Range("M136").Value = -4523.3334
$ws.Range("N136").Value = -48600


# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 12
$ws.Range("H12").Value = 4300
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 4300
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 4300
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -4636

# Row 86
$ws.Range("H86").Value = 439193.75
$ws.Range("I86").Value = 1578.5
$ws.Range("J86").Value = 3502500.5
$ws.Range("K86").Value = 1578.5
$ws.Range("L86").Value = 3502500.5
$ws.Range("M86").Value = -455.5
$ws.Range("N86").Value = -3504746.5

# Row 89
$ws.Range("H89").Value = 439193.75
$ws.Range("I89").Value = 1578.5
$ws.Range("J89").Value = 3502500.5
$ws.Range("K89").Value = 7892.5
$ws.Range("L89").Value = 17512502.5
$ws.Range("M89").Value = -2276.5
$ws.Range("N89").Value = -17523734.5

# Row 134
$ws.Range("H134").Value = 17896574
$ws.Range("I134").Value = 29413704
$ws.Range("J134").Value = 97375.27
$ws.Range("K134").Value = 88241112
$ws.Range("L134").Value = 292125.81
$ws.Range("M134").Value = -88238577
$ws.Range("N134").Value = -297195.81


# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 8068.1113
$ws.Range("I58").Value = 1435.6666
$ws.Range("J58").Value = 21333
$ws.Range("K58").Value = 1435.6666
$ws.Range("L58").Value = 21333
$ws.Range("M58").Value = -1232.6666
$ws.Range("N58").Value = -21739

# Row 93
$ws.Range("H93").Value = 12333.167
$ws.Range("I93").Value = 4799.8
$ws.Range("J93").Value = 50000
$ws.Range("K93").Value = 4799.8
$ws.Range("L93").Value = 50000
$ws.Range("M93").Value = -2927.8
$ws.Range("N93").Value = -53744

# Row 99
$ws.Range("H99").Value = 2855.6667
$ws.Range("I99").Value = 2775.3333
$ws.Range("J99").Value = 3016.3333
$ws.Range("K99").Value = 2775.3333
$ws.Range("L99").Value = 3016.3333
$ws.Range("M99").Value = -1277.3333
$ws.Range("N99").Value = -6012.3333

# Row 126
$ws.Range("H126").Value = 2855.6667
$ws.Range("I126").Value = 2775.3333
$ws.Range("J126").Value = 3016.3333
$ws.Range("K126").Value = 8325.999899999999
$ws.Range("L126").Value = 9048.999899999999
$ws.Range("M126").Value = -5855.999899999999
$ws.Range("N126").Value = -13988.9999

# Row 134
$ws.Range("H134").Value = 27779766
$ws.Range("I134").Value = 1814.2
$ws.Range("J134").Value = 62502210
$ws.Range("K134").Value = 5442.6
$ws.Range("L134").Value = 187506630
$ws.Range("M134").Value = -2907.6
$ws.Range("N134").Value = -187511700

# Row 136
$ws.Range("H136").Value = 8068.1113
$ws.Range("I136").Value = 1435.6666
$ws.Range("J136").Value = 21333
$ws.Range("K136").Value = 4306.9998
$ws.Range("L136").Value = 63999
$ws.Range("M136").Value = -1756.9998
$ws.Range("N136").Value = -69099


# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 91
$ws.Range("I38").Value = 80
$ws.Range("J38").Value = 102
$ws.Range("K38").Value = 240
$ws.Range("L38").Value = 306
$ws.Range("M38").Value = 107
$ws.Range("N38").Value = -1000

# Row 63
$ws.Range("H63").Value = 2761.3333
$ws.Range("I63").Value = 1512
$ws.Range("J63").Value = 3011.2
$ws.Range("K63").Value = 4536
$ws.Range("L63").Value = 9033.599999999999
$ws.Range("M63").Value = -3787
$ws.Range("N63").Value = -10531.6

# Row 66
$ws.Range("H66").Value = 2761.3333
$ws.Range("I66").Value = 1512
$ws.Range("J66").Value = 3011.2
$ws.Range("K66").Value = 13608
$ws.Range("L66").Value = 27100.8
$ws.Range("M66").Value = -9864
$ws.Range("N66").Value = -34588.8

# Row 107
$ws.Range("H107").Value = 655.38464
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 752.5
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 2257.5
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -6097.5


# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 13
$ws.Range("H13").Value = 8086.857
$ws.Range("I13").Value = 121.6
$ws.Range("J13").Value = 28000
$ws.Range("K13").Value = 121.6
$ws.Range("L13").Value = 28000
$ws.Range("M13").Value = 17.40000000000001
$ws.Range("N13").Value = -28278

# Row 70
$ws.Range("H70").Value = 4240.5625
$ws.Range("I70").Value = 4189.3687
$ws.Range("K70").Value = 4189.3687
$ws.Range("M70").Value = -3919.3687

# Row 73
$ws.Range("H73").Value = 4240.5625
$ws.Range("I73").Value = 4189.3687
$ws.Range("K73").Value = 4189.3687
$ws.Range("M73").Value = -3253.3687

# Row 113
$ws.Range("H113").Value = 1287.8572
$ws.Range("I113").Value = 1317
$ws.Range("J113").Value = 1113
$ws.Range("K113").Value = 1317
$ws.Range("L113").Value = 1113
$ws.Range("M113").Value = 853
$ws.Range("N113").Value = -5453


# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 3250
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 5500
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 5500
$ws.Range("M16").Value = -830
$ws.Range("N16").Value = -5840

# Row 48
$ws.Range("H48").Value = 51627.6
$ws.Range("J48").Value = 51627.6
$ws.Range("L48").Value = 51627.6
$ws.Range("N48").Value = -52949.6

# Row 122
$ws.Range("H122").Value = 3031.6
$ws.Range("I122").Value = 2810.4
$ws.Range("K122").Value = 8431.200000000001
$ws.Range("M122").Value = -5981.200000000001

# Row 136
$ws.Range("H136").Value = 557988.6
$ws.Range("I136").Value = 1251837.5
$ws.Range("K136").Value = 3755512.5
$ws.Range("M136").Value = -3752962.5


# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 3702702.8
$ws.Range("I136").Value = 3969802.8
$ws.Range("J136").Value = 2500752.5
$ws.Range("K136").Value = 11909408.4
$ws.Range("L136").Value = 7502257.5
$ws.Range("M136").Value = -11906858.4
$ws.Range("N136").Value = -7507357.5

